# Apply a channel-flip to the grids: swap D-column values pair-wise
# (rows 3-14, corresponding to grid index 0-11) since the grids were
# physically swapped pair-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairwise swap the "Channel" values for grids (column D) in rows 3..14
for ($r = 3; $r -le 13; $r += 2) {
    $top = $ws.Cells.Item($r, 4)
    $bottom = $ws.Cells.Item($r + 1, 4)
    $topValue = $top.Value()
    $bottomValue = $bottom.Value()
    $top.Value = $bottomValue
    $bottom.Value = $topValue
}

# Update the view: move the selection to D15
$ws.Range("D15").Select()
